$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap home/away match data between paired rows (F:V), sorted differently after refresh ---
# Row 36
$ws.Cells.Item(36, 6).Value = 'PSV'
$ws.Cells.Item(36, 7).Value = 4
$ws.Cells.Item(36, 8).Value = 'Nijmegen'
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 1.24
$ws.Cells.Item(36, 11).Value = '04/09/2023 08:43'
$ws.Cells.Item(36, 12).Value = 1.15
$ws.Cells.Item(36, 13).Value = '16/09/2023 19:38'
$ws.Cells.Item(36, 14).Value = 7.01
$ws.Cells.Item(36, 15).Value = '04/09/2023 08:43'
$ws.Cells.Item(36, 16).Value = 9.119999999999999
$ws.Cells.Item(36, 17).Value = '16/09/2023 19:57'
$ws.Cells.Item(36, 18).Value = 11.5
$ws.Cells.Item(36, 19).Value = '04/09/2023 08:43'
$ws.Cells.Item(36, 20).Value = 17.78
$ws.Cells.Item(36, 21).Value = '16/09/2023 19:57'
$ws.Cells.Item(36, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-nijmegen/8zF9rNhq/'

# Row 37
$ws.Cells.Item(37, 6).Value = 'Sittard'
$ws.Cells.Item(37, 7).Value = 3
$ws.Cells.Item(37, 8).Value = 'FC Volendam'
$ws.Cells.Item(37, 9).Value = 1
$ws.Cells.Item(37, 10).Value = 1.67
$ws.Cells.Item(37, 11).Value = '04/09/2023 08:43'
$ws.Cells.Item(37, 12).Value = 1.47
$ws.Cells.Item(37, 13).Value = '16/09/2023 19:59'
$ws.Cells.Item(37, 14).Value = 4.44
$ws.Cells.Item(37, 15).Value = '04/09/2023 08:43'
$ws.Cells.Item(37, 16).Value = 4.96
$ws.Cells.Item(37, 17).Value = '16/09/2023 19:59'
$ws.Cells.Item(37, 18).Value = 4.56
$ws.Cells.Item(37, 19).Value = '04/09/2023 08:43'
$ws.Cells.Item(37, 20).Value = 6.75
$ws.Cells.Item(37, 21).Value = '16/09/2023 19:59'
$ws.Cells.Item(37, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sittard-fc-volendam/StXtm3pS/'

# Row 40
$ws.Cells.Item(40, 6).Value = 'Twente'
$ws.Cells.Item(40, 7).Value = 3
$ws.Cells.Item(40, 8).Value = 'Ajax'
$ws.Cells.Item(40, 9).Value = 1
$ws.Cells.Item(40, 10).Value = 3.43
$ws.Cells.Item(40, 11).Value = '04/09/2023 08:43'
$ws.Cells.Item(40, 12).Value = 2.47
$ws.Cells.Item(40, 13).Value = '17/09/2023 14:29'
$ws.Cells.Item(40, 14).Value = 4.07
$ws.Cells.Item(40, 15).Value = '04/09/2023 08:43'
$ws.Cells.Item(40, 16).Value = 3.6
$ws.Cells.Item(40, 17).Value = '17/09/2023 14:28'
$ws.Cells.Item(40, 18).Value = 2.02
$ws.Cells.Item(40, 19).Value = '04/09/2023 08:43'
$ws.Cells.Item(40, 20).Value = 2.9
$ws.Cells.Item(40, 21).Value = '17/09/2023 14:29'
$ws.Cells.Item(40, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-ajax/dde444F2/'

# Row 41
$ws.Cells.Item(41, 6).Value = 'Excelsior'
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 'Almere City'
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 1.93
$ws.Cells.Item(41, 11).Value = '04/09/2023 08:43'
$ws.Cells.Item(41, 12).Value = 2.04
$ws.Cells.Item(41, 13).Value = '17/09/2023 14:29'
$ws.Cells.Item(41, 14).Value = 3.96
$ws.Cells.Item(41, 15).Value = '04/09/2023 08:43'
$ws.Cells.Item(41, 16).Value = 3.95
$ws.Cells.Item(41, 17).Value = '17/09/2023 14:28'
$ws.Cells.Item(41, 18).Value = 3.8
$ws.Cells.Item(41, 19).Value = '04/09/2023 08:43'
$ws.Cells.Item(41, 20).Value = 3.49
$ws.Cells.Item(41, 21).Value = '17/09/2023 14:29'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-almere-city/fiHLuLx2/'

# Row 56
$ws.Cells.Item(56, 6).Value = 'Utrecht'
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 'Almere City'
$ws.Cells.Item(56, 9).Value = 2
$ws.Cells.Item(56, 10).Value = 1.51
$ws.Cells.Item(56, 11).Value = '23/09/2023 19:12'
$ws.Cells.Item(56, 12).Value = 1.57
$ws.Cells.Item(56, 13).Value = '30/09/2023 18:44'
$ws.Cells.Item(56, 14).Value = 4.75
$ws.Cells.Item(56, 15).Value = '23/09/2023 19:12'
$ws.Cells.Item(56, 16).Value = 4.4
$ws.Cells.Item(56, 17).Value = '30/09/2023 18:44'
$ws.Cells.Item(56, 18).Value = 5.79
$ws.Cells.Item(56, 19).Value = '23/09/2023 19:12'
$ws.Cells.Item(56, 20).Value = 5.93
$ws.Cells.Item(56, 21).Value = '30/09/2023 18:44'
$ws.Cells.Item(56, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/utrecht-almere-city/dv2Y7wMf/'

# Row 57
$ws.Cells.Item(57, 6).Value = 'PSV'
$ws.Cells.Item(57, 7).Value = 3
$ws.Cells.Item(57, 8).Value = 'FC Volendam'
$ws.Cells.Item(57, 9).Value = 1
$ws.Cells.Item(57, 10).Value = 1.06
$ws.Cells.Item(57, 11).Value = '27/09/2023 18:12'
$ws.Cells.Item(57, 12).Value = 1.04
$ws.Cells.Item(57, 13).Value = '30/09/2023 14:21'
$ws.Cells.Item(57, 14).Value = 17.77
$ws.Cells.Item(57, 15).Value = '27/09/2023 18:12'
$ws.Cells.Item(57, 16).Value = 25.88
$ws.Cells.Item(57, 17).Value = '30/09/2023 18:44'
$ws.Cells.Item(57, 18).Value = 24.06
$ws.Cells.Item(57, 19).Value = '27/09/2023 18:12'
$ws.Cells.Item(57, 20).Value = 42.53
$ws.Cells.Item(57, 21).Value = '30/09/2023 18:44'
$ws.Cells.Item(57, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-fc-volendam/EFCT8J6l/'

# Row 60
$ws.Cells.Item(60, 6).Value = 'Heracles'
$ws.Cells.Item(60, 7).Value = 2
$ws.Cells.Item(60, 8).Value = 'Zwolle'
$ws.Cells.Item(60, 9).Value = 1
$ws.Cells.Item(60, 10).Value = 2.09
$ws.Cells.Item(60, 11).Value = '28/09/2023 19:12'
$ws.Cells.Item(60, 12).Value = 2.44
$ws.Cells.Item(60, 13).Value = '01/10/2023 14:24'
$ws.Cells.Item(60, 14).Value = 4.05
$ws.Cells.Item(60, 15).Value = '28/09/2023 19:12'
$ws.Cells.Item(60, 16).Value = 3.69
$ws.Cells.Item(60, 17).Value = '01/10/2023 14:27'
$ws.Cells.Item(60, 18).Value = 3.19
$ws.Cells.Item(60, 19).Value = '28/09/2023 19:12'
$ws.Cells.Item(60, 20).Value = 2.88
$ws.Cells.Item(60, 21).Value = '01/10/2023 14:27'
$ws.Cells.Item(60, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/heracles-zwolle/rmALTKrE/'

# Row 61
$ws.Cells.Item(61, 6).Value = 'Excelsior'
$ws.Cells.Item(61, 7).Value = 2
$ws.Cells.Item(61, 8).Value = 'Sparta Rotterdam'
$ws.Cells.Item(61, 9).Value = 1
$ws.Cells.Item(61, 10).Value = 2.83
$ws.Cells.Item(61, 11).Value = '24/09/2023 13:42'
$ws.Cells.Item(61, 12).Value = 3.17
$ws.Cells.Item(61, 13).Value = '01/10/2023 14:20'
$ws.Cells.Item(61, 14).Value = 3.52
$ws.Cells.Item(61, 15).Value = '24/09/2023 13:42'
$ws.Cells.Item(61, 16).Value = 3.54
$ws.Cells.Item(61, 17).Value = '01/10/2023 14:20'
$ws.Cells.Item(61, 18).Value = 2.49
$ws.Cells.Item(61, 19).Value = '24/09/2023 13:42'
$ws.Cells.Item(61, 20).Value = 2.33
$ws.Cells.Item(61, 21).Value = '01/10/2023 14:20'
$ws.Cells.Item(61, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-sparta-rotterdam/KCXlChDt/'

# Row 69
$ws.Cells.Item(69, 6).Value = 'Sittard'
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 'Twente'
$ws.Cells.Item(69, 9).Value = 3
$ws.Cells.Item(69, 10).Value = 4.64
$ws.Cells.Item(69, 11).Value = '01/10/2023 16:12'
$ws.Cells.Item(69, 12).Value = 4.8
$ws.Cells.Item(69, 13).Value = '08/10/2023 14:29'
$ws.Cells.Item(69, 14).Value = 4.14
$ws.Cells.Item(69, 15).Value = '01/10/2023 16:12'
$ws.Cells.Item(69, 16).Value = 4.04
$ws.Cells.Item(69, 17).Value = '08/10/2023 14:29'
$ws.Cells.Item(69, 18).Value = 1.71
$ws.Cells.Item(69, 19).Value = '01/10/2023 16:12'
$ws.Cells.Item(69, 20).Value = 1.74
$ws.Cells.Item(69, 21).Value = '08/10/2023 14:29'
$ws.Cells.Item(69, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sittard-twente/SxcSPBst/'

# Row 70
$ws.Cells.Item(70, 6).Value = 'Ajax'
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = 'AZ Alkmaar'
$ws.Cells.Item(70, 9).Value = 2
$ws.Cells.Item(70, 10).Value = 2.24
$ws.Cells.Item(70, 11).Value = '01/10/2023 16:12'
$ws.Cells.Item(70, 12).Value = 2.7
$ws.Cells.Item(70, 13).Value = '08/10/2023 14:29'
$ws.Cells.Item(70, 14).Value = 3.77
$ws.Cells.Item(70, 15).Value = '01/10/2023 16:12'
$ws.Cells.Item(70, 16).Value = 3.74
$ws.Cells.Item(70, 17).Value = '08/10/2023 14:29'
$ws.Cells.Item(70, 18).Value = 3.06
$ws.Cells.Item(70, 19).Value = '01/10/2023 16:12'
$ws.Cells.Item(70, 20).Value = 2.56
$ws.Cells.Item(70, 21).Value = '08/10/2023 14:29'
$ws.Cells.Item(70, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/ajax-az-alkmaar/jHLM4SBU/'

# Row 72
$ws.Cells.Item(72, 6).Value = 'PSV'
$ws.Cells.Item(72, 7).Value = 3
$ws.Cells.Item(72, 8).Value = 'Sittard'
$ws.Cells.Item(72, 9).Value = 1
$ws.Cells.Item(72, 10).Value = 1.1
$ws.Cells.Item(72, 11).Value = '09/10/2023 14:42'
$ws.Cells.Item(72, 12).Value = 1.14
$ws.Cells.Item(72, 13).Value = '21/10/2023 18:02'
$ws.Cells.Item(72, 14).Value = 12.21
$ws.Cells.Item(72, 15).Value = '09/10/2023 14:42'
$ws.Cells.Item(72, 16).Value = 9.73
$ws.Cells.Item(72, 17).Value = '21/10/2023 18:41'
$ws.Cells.Item(72, 18).Value = 18.64
$ws.Cells.Item(72, 19).Value = '09/10/2023 14:42'
$ws.Cells.Item(72, 20).Value = 17.33
$ws.Cells.Item(72, 21).Value = '21/10/2023 18:41'
$ws.Cells.Item(72, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/psv-sittard/2ekuN9Rb/'

# Row 73
$ws.Cells.Item(73, 6).Value = 'Waalwijk'
$ws.Cells.Item(73, 7).Value = 2
$ws.Cells.Item(73, 8).Value = 'FC Volendam'
$ws.Cells.Item(73, 9).Value = 1
$ws.Cells.Item(73, 10).Value = 1.5
$ws.Cells.Item(73, 11).Value = '09/10/2023 14:42'
$ws.Cells.Item(73, 12).Value = 1.76
$ws.Cells.Item(73, 13).Value = '21/10/2023 18:44'
$ws.Cells.Item(73, 14).Value = 4.82
$ws.Cells.Item(73, 15).Value = '09/10/2023 14:42'
$ws.Cells.Item(73, 16).Value = 4.24
$ws.Cells.Item(73, 17).Value = '21/10/2023 18:41'
$ws.Cells.Item(73, 18).Value = 5.87
$ws.Cells.Item(73, 19).Value = '09/10/2023 14:42'
$ws.Cells.Item(73, 20).Value = 4.41
$ws.Cells.Item(73, 21).Value = '21/10/2023 18:44'
$ws.Cells.Item(73, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/waalwijk-fc-volendam/dIlqMTt5/'

# Row 95
$ws.Cells.Item(95, 6).Value = 'Nijmegen'
$ws.Cells.Item(95, 7).Value = 3
$ws.Cells.Item(95, 8).Value = 'FC Volendam'
$ws.Cells.Item(95, 9).Value = 3
$ws.Cells.Item(95, 10).Value = 1.6
$ws.Cells.Item(95, 11).Value = '02/11/2023 20:12'
$ws.Cells.Item(95, 12).Value = 1.56
$ws.Cells.Item(95, 13).Value = '05/11/2023 14:27'
$ws.Cells.Item(95, 14).Value = 4.52
$ws.Cells.Item(95, 15).Value = '02/11/2023 20:12'
$ws.Cells.Item(95, 16).Value = 4.71
$ws.Cells.Item(95, 17).Value = '05/11/2023 14:27'
$ws.Cells.Item(95, 18).Value = 5.11
$ws.Cells.Item(95, 19).Value = '02/11/2023 20:12'
$ws.Cells.Item(95, 20).Value = 5.55
$ws.Cells.Item(95, 21).Value = '05/11/2023 14:29'
$ws.Cells.Item(95, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/nijmegen-fc-volendam/vyc9qOvT/'

# Row 96
$ws.Cells.Item(96, 6).Value = 'Zwolle'
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 'Sittard'
$ws.Cells.Item(96, 9).Value = 0
$ws.Cells.Item(96, 10).Value = 2.08
$ws.Cells.Item(96, 11).Value = '28/10/2023 21:12'
$ws.Cells.Item(96, 12).Value = 2.31
$ws.Cells.Item(96, 13).Value = '05/11/2023 14:21'
$ws.Cells.Item(96, 14).Value = 3.65
$ws.Cells.Item(96, 15).Value = '28/10/2023 21:12'
$ws.Cells.Item(96, 16).Value = 3.5
$ws.Cells.Item(96, 17).Value = '05/11/2023 14:28'
$ws.Cells.Item(96, 18).Value = 3.51
$ws.Cells.Item(96, 19).Value = '28/10/2023 21:12'
$ws.Cells.Item(96, 20).Value = 3.24
$ws.Cells.Item(96, 21).Value = '05/11/2023 14:28'
$ws.Cells.Item(96, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-sittard/tpqyurmp/'

# --- Append 9 new match rows (108-116), copying formats from row 107 first ---
$ws.Range("A107:V107").Copy() | Out-Null
$ws.Range("A108:V116").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 108
$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = 'netherlands'
$ws.Cells.Item(108, 3).Value = 'eredivisie'
$ws.Cells.Item(108, 4).Value = '2023-2024'
$ws.Cells.Item(108, 5).Value = 45255.6875
$ws.Cells.Item(108, 6).Value = 'Excelsior'
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 'Feyenoord'
$ws.Cells.Item(108, 9).Value = 4
$ws.Cells.Item(108, 10).Value = 8.35
$ws.Cells.Item(108, 11).Value = '12/11/2023 20:12'
$ws.Cells.Item(108, 12).Value = 14.12
$ws.Cells.Item(108, 13).Value = '25/11/2023 16:24'
$ws.Cells.Item(108, 14).Value = 6.28
$ws.Cells.Item(108, 15).Value = '12/11/2023 20:12'
$ws.Cells.Item(108, 16).Value = 8.289999999999999
$ws.Cells.Item(108, 17).Value = '25/11/2023 16:24'
$ws.Cells.Item(108, 18).Value = 1.27
$ws.Cells.Item(108, 19).Value = '12/11/2023 20:12'
$ws.Cells.Item(108, 20).Value = 1.18
$ws.Cells.Item(108, 21).Value = '25/11/2023 16:21'
$ws.Cells.Item(108, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/excelsior-feyenoord/8WvMN1Z6/'

# Row 109
$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = 'netherlands'
$ws.Cells.Item(109, 3).Value = 'eredivisie'
$ws.Cells.Item(109, 4).Value = '2023-2024'
$ws.Cells.Item(109, 5).Value = 45255.78125
$ws.Cells.Item(109, 6).Value = 'Twente'
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 'PSV'
$ws.Cells.Item(109, 9).Value = 3
$ws.Cells.Item(109, 10).Value = 3.37
$ws.Cells.Item(109, 11).Value = '12/11/2023 15:42'
$ws.Cells.Item(109, 12).Value = 3.65
$ws.Cells.Item(109, 13).Value = '25/11/2023 18:41'
$ws.Cells.Item(109, 14).Value = 4.11
$ws.Cells.Item(109, 15).Value = '12/11/2023 15:42'
$ws.Cells.Item(109, 16).Value = 3.76
$ws.Cells.Item(109, 17).Value = '25/11/2023 18:41'
$ws.Cells.Item(109, 18).Value = 2
$ws.Cells.Item(109, 19).Value = '12/11/2023 15:42'
$ws.Cells.Item(109, 20).Value = 2.05
$ws.Cells.Item(109, 21).Value = '25/11/2023 18:41'
$ws.Cells.Item(109, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/twente-psv/0QR7o5RD/'

# Row 110
$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = 'netherlands'
$ws.Cells.Item(110, 3).Value = 'eredivisie'
$ws.Cells.Item(110, 4).Value = '2023-2024'
$ws.Cells.Item(110, 5).Value = 45255.78125
$ws.Cells.Item(110, 6).Value = 'Heerenveen'
$ws.Cells.Item(110, 7).Value = 3
$ws.Cells.Item(110, 8).Value = 'Sittard'
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 1.92
$ws.Cells.Item(110, 11).Value = '12/11/2023 15:42'
$ws.Cells.Item(110, 12).Value = 2.02
$ws.Cells.Item(110, 13).Value = '25/11/2023 18:42'
$ws.Cells.Item(110, 14).Value = 3.86
$ws.Cells.Item(110, 15).Value = '12/11/2023 15:42'
$ws.Cells.Item(110, 16).Value = 3.44
$ws.Cells.Item(110, 17).Value = '25/11/2023 18:42'
$ws.Cells.Item(110, 18).Value = 3.85
$ws.Cells.Item(110, 19).Value = '12/11/2023 15:42'
$ws.Cells.Item(110, 20).Value = 4.13
$ws.Cells.Item(110, 21).Value = '25/11/2023 18:42'
$ws.Cells.Item(110, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/heerenveen-sittard/nNwQMLlD/'

# Row 111
$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = 'netherlands'
$ws.Cells.Item(111, 3).Value = 'eredivisie'
$ws.Cells.Item(111, 4).Value = '2023-2024'
$ws.Cells.Item(111, 5).Value = 45255.875
$ws.Cells.Item(111, 6).Value = 'Ajax'
$ws.Cells.Item(111, 7).Value = 5
$ws.Cells.Item(111, 8).Value = 'Vitesse'
$ws.Cells.Item(111, 9).Value = 0
$ws.Cells.Item(111, 10).Value = 1.37
$ws.Cells.Item(111, 11).Value = '13/11/2023 09:49'
$ws.Cells.Item(111, 12).Value = 1.31
$ws.Cells.Item(111, 13).Value = '25/11/2023 20:54'
$ws.Cells.Item(111, 14).Value = 5.63
$ws.Cells.Item(111, 15).Value = '13/11/2023 09:49'
$ws.Cells.Item(111, 16).Value = 6.46
$ws.Cells.Item(111, 17).Value = '25/11/2023 20:38'
$ws.Cells.Item(111, 18).Value = 6.3
$ws.Cells.Item(111, 19).Value = '13/11/2023 09:49'
$ws.Cells.Item(111, 20).Value = 8.4
$ws.Cells.Item(111, 21).Value = '25/11/2023 20:38'
$ws.Cells.Item(111, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/ajax-vitesse/ETVBpPtK/'

# Row 112
$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = 'netherlands'
$ws.Cells.Item(112, 3).Value = 'eredivisie'
$ws.Cells.Item(112, 4).Value = '2023-2024'
$ws.Cells.Item(112, 5).Value = 45255.875
$ws.Cells.Item(112, 6).Value = 'Zwolle'
$ws.Cells.Item(112, 7).Value = 1
$ws.Cells.Item(112, 8).Value = 'Waalwijk'
$ws.Cells.Item(112, 9).Value = 2
$ws.Cells.Item(112, 10).Value = 1.78
$ws.Cells.Item(112, 11).Value = '12/11/2023 15:42'
$ws.Cells.Item(112, 12).Value = 2.04
$ws.Cells.Item(112, 13).Value = '25/11/2023 20:37'
$ws.Cells.Item(112, 14).Value = 3.97
$ws.Cells.Item(112, 15).Value = '12/11/2023 15:42'
$ws.Cells.Item(112, 16).Value = 3.78
$ws.Cells.Item(112, 17).Value = '25/11/2023 20:37'
$ws.Cells.Item(112, 18).Value = 4.39
$ws.Cells.Item(112, 19).Value = '12/11/2023 15:42'
$ws.Cells.Item(112, 20).Value = 3.65
$ws.Cells.Item(112, 21).Value = '25/11/2023 20:37'
$ws.Cells.Item(112, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/zwolle-waalwijk/dhLGqqdQ/'

# Row 113
$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = 'netherlands'
$ws.Cells.Item(113, 3).Value = 'eredivisie'
$ws.Cells.Item(113, 4).Value = '2023-2024'
$ws.Cells.Item(113, 5).Value = 45256.51041666666
$ws.Cells.Item(113, 6).Value = 'Almere City'
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 'Heracles'
$ws.Cells.Item(113, 9).Value = 5
$ws.Cells.Item(113, 10).Value = 2.11
$ws.Cells.Item(113, 11).Value = '12/11/2023 15:42'
$ws.Cells.Item(113, 12).Value = 1.9
$ws.Cells.Item(113, 13).Value = '26/11/2023 12:14'
$ws.Cells.Item(113, 14).Value = 3.74
$ws.Cells.Item(113, 15).Value = '12/11/2023 15:42'
$ws.Cells.Item(113, 16).Value = 3.92
$ws.Cells.Item(113, 17).Value = '26/11/2023 12:14'
$ws.Cells.Item(113, 18).Value = 3.37
$ws.Cells.Item(113, 19).Value = '12/11/2023 15:42'
$ws.Cells.Item(113, 20).Value = 4.04
$ws.Cells.Item(113, 21).Value = '26/11/2023 12:14'
$ws.Cells.Item(113, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-heracles/b17ru34s/'

# Row 114
$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = 'netherlands'
$ws.Cells.Item(114, 3).Value = 'eredivisie'
$ws.Cells.Item(114, 4).Value = '2023-2024'
$ws.Cells.Item(114, 5).Value = 45256.60416666666
$ws.Cells.Item(114, 6).Value = 'AZ Alkmaar'
$ws.Cells.Item(114, 7).Value = 3
$ws.Cells.Item(114, 8).Value = 'FC Volendam'
$ws.Cells.Item(114, 9).Value = 0
$ws.Cells.Item(114, 10).Value = 1.16
$ws.Cells.Item(114, 11).Value = '12/11/2023 17:12'
$ws.Cells.Item(114, 12).Value = 1.18
$ws.Cells.Item(114, 13).Value = '26/11/2023 14:26'
$ws.Cells.Item(114, 14).Value = 8.99
$ws.Cells.Item(114, 15).Value = '12/11/2023 17:12'
$ws.Cells.Item(114, 16).Value = 8.44
$ws.Cells.Item(114, 17).Value = '26/11/2023 14:29'
$ws.Cells.Item(114, 18).Value = 14.18
$ws.Cells.Item(114, 19).Value = '12/11/2023 17:12'
$ws.Cells.Item(114, 20).Value = 14.91
$ws.Cells.Item(114, 21).Value = '26/11/2023 14:29'
$ws.Cells.Item(114, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/az-alkmaar-fc-volendam/pI5jwsZg/'

# Row 115
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = 'netherlands'
$ws.Cells.Item(115, 3).Value = 'eredivisie'
$ws.Cells.Item(115, 4).Value = '2023-2024'
$ws.Cells.Item(115, 5).Value = 45256.60416666666
$ws.Cells.Item(115, 6).Value = 'Nijmegen'
$ws.Cells.Item(115, 7).Value = 1
$ws.Cells.Item(115, 8).Value = 'G.A. Eagles'
$ws.Cells.Item(115, 9).Value = 1
$ws.Cells.Item(115, 10).Value = 2.08
$ws.Cells.Item(115, 11).Value = '12/11/2023 15:42'
$ws.Cells.Item(115, 12).Value = 2.25
$ws.Cells.Item(115, 13).Value = '26/11/2023 14:29'
$ws.Cells.Item(115, 14).Value = 3.85
$ws.Cells.Item(115, 15).Value = '12/11/2023 15:42'
$ws.Cells.Item(115, 16).Value = 3.67
$ws.Cells.Item(115, 17).Value = '26/11/2023 14:25'
$ws.Cells.Item(115, 18).Value = 3.34
$ws.Cells.Item(115, 19).Value = '12/11/2023 15:42'
$ws.Cells.Item(115, 20).Value = 3.22
$ws.Cells.Item(115, 21).Value = '26/11/2023 14:29'
$ws.Cells.Item(115, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/nijmegen-g-a-eagles/6R6nvNJm/'

# Row 116
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = 'netherlands'
$ws.Cells.Item(116, 3).Value = 'eredivisie'
$ws.Cells.Item(116, 4).Value = '2023-2024'
$ws.Cells.Item(116, 5).Value = 45256.69791666666
$ws.Cells.Item(116, 6).Value = 'Sparta Rotterdam'
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 'Utrecht'
$ws.Cells.Item(116, 9).Value = 2
$ws.Cells.Item(116, 10).Value = 2.34
$ws.Cells.Item(116, 11).Value = '12/11/2023 20:12'
$ws.Cells.Item(116, 12).Value = 2.52
$ws.Cells.Item(116, 13).Value = '26/11/2023 16:41'
$ws.Cells.Item(116, 14).Value = 3.31
$ws.Cells.Item(116, 15).Value = '12/11/2023 20:12'
$ws.Cells.Item(116, 16).Value = 3.38
$ws.Cells.Item(116, 17).Value = '26/11/2023 16:40'
$ws.Cells.Item(116, 18).Value = 3.07
$ws.Cells.Item(116, 19).Value = '12/11/2023 20:12'
$ws.Cells.Item(116, 20).Value = 2.98
$ws.Cells.Item(116, 21).Value = '26/11/2023 16:41'
$ws.Cells.Item(116, 22).Value = 'https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-utrecht/Steex1l0/'
